$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$text) {
    # Assigning a quoted-string formula and then collapsing it to a value via
    # PasteSpecial(xlPasteValues) keeps the cell's original style untouched
    # and avoids Excel's automatic text->number coercion that a direct
    # .Value assignment would trigger for numeric-looking strings
    # (e.g. "1.00" -> 1, "0.650" -> 0.65).
    $cell = $ws.Range($cellRef)
    $cell.Formula = "=""" + $text + """"
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

Set-TextValue "D2" "60.828.57"
$ws.Range("E2").Value = "  -3.79%  "
Set-TextValue "D3" "2.928.22"
$ws.Range("E3").Value = "  -3.04%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue "D5" "530.29"
$ws.Range("E5").Value = "  -4.77%  "
Set-TextValue "D6" "145.86"
$ws.Range("E6").Value = "  -5.62%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.01%  "
Set-TextValue "D9" "2.936.32"
$ws.Range("E9").Value = "  -3.03%  "
$ws.Range("E10").Value = "  -3.38%  "
$ws.Range("E11").Value = "  -7.46%  "
Set-TextValue "D12" "0.356"
$ws.Range("E12").Value = "  -2.93%  "
Set-TextValue "D13" "3.436.52"
$ws.Range("E13").Value = "  -3.27%  "
$ws.Range("E14").Value = "  +1.69%  "
Set-TextValue "D15" "60.944.63"
$ws.Range("E15").Value = "  -3.67%  "
Set-TextValue "D16" "23.03"
$ws.Range("E16").Value = "  -4.44%  "
Set-TextValue "D17" "2.924.90"
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("E18").Value = "  -5.94%  "
$ws.Range("E19").Value = "  -1.74%  "
Set-TextValue "D20" "11.71"
$ws.Range("E20").Value = "  -2.41%  "
Set-TextValue "D21" "364.72"
$ws.Range("E21").Value = "  -8.55%  "
Set-TextValue "D22" "6.52"
$ws.Range("E22").Value = "  -2.16%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -1.92%  "
Set-TextValue "D25" "64.86"
$ws.Range("E25").Value = "  -0.68%  "
Set-TextValue "D26" "3.049.59"
$ws.Range("E26").Value = "  -3.46%  "
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("E28").Value = "  -2.60%  "
Set-TextValue "D29" "0.998"
$ws.Range("E29").Value = "  -0.13%  "
Set-TextValue "D30" "0.0₃0875"
$ws.Range("E30").Value = "  -11.29%  "
Set-TextValue "D31" "7.75"
$ws.Range("E31").Value = "  -10.94%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  -4.84%  "
Set-TextValue "D34" "19.89"
$ws.Range("E34").Value = "  -2.52%  "
Set-TextValue "D35" "156.03"
$ws.Range("E35").Value = "  -3.30%  "
Set-TextValue "D36" "4.42"
$ws.Range("E36").Value = "  -6.70%  "
Set-TextValue "D37" "5.67"
$ws.Range("E37").Value = "  -6.07%  "
$ws.Range("E38").Value = "  -9.20%  "
$ws.Range("E39").Value = "  -6.63%  "
Set-TextValue "D40" "38.27"
$ws.Range("E40").Value = "  +1.34%  "
$ws.Range("E41").Value = "  -5.77%  "
Set-TextValue "D42" "2.358.85"
$ws.Range("E42").Value = "  -7.45%  "
Set-TextValue "D43" "3.74"
$ws.Range("E43").Value = "  -5.36%  "
Set-TextValue "D44" "0.650"
$ws.Range("E44").Value = "  -2.92%  "
Set-TextValue "D45" "21.09"
$ws.Range("E45").Value = "  -7.95%  "
Set-TextValue "D46" "0.0574"
$ws.Range("E46").Value = "  -4.16%  "
Set-TextValue "D47" "1.00"
$ws.Range("E47").Value = "  +0.23%  "
Set-TextValue "D48" "4.97"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("E49").Value = "  -6.69%  "
Set-TextValue "D50" "10.36"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("E51").Value = "  -1.61%  "

$excel.CutCopyMode = $false
